$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.230.51'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.553.79'
$ws.Range('E3').Value = '  +3.52%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.63'
$ws.Range('E5').Value = '  +0.62%  '
$ws.Range('E6').Value = '  +3.49%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.588'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.552.28'
$ws.Range('E9').Value = '  +3.53%  '
$ws.Range('E10').Value = '  +0.82%  '
$ws.Range('E11').Value = '  -1.68%  '
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.62'
$ws.Range('E14').Value = '  +3.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.010.09'
$ws.Range('E15').Value = '  +3.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.163.99'
$ws.Range('E16').Value = '  +0.59%  '
$ws.Range('E17').Value = '  +2.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.548.62'
$ws.Range('E18').Value = '  +3.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.46'
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '335.82'
$ws.Range('E20').Value = '  -1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.33'
$ws.Range('E21').Value = '  +1.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.78'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.29'
$ws.Range('E24').Value = '  -0.43%  '
$ws.Range('E25').Value = '  +9.88%  '
$ws.Range('E26').Value = '  -1.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.49'
$ws.Range('E27').Value = '  +7.48%  '
$ws.Range('E28').Value = '  +5.52%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  +8.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0823'
$ws.Range('E31').Value = '  +2.95%  '
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '176.78'
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('E34').Value = '  +4.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '413.05'
$ws.Range('E35').Value = '  +12.51%  '
$ws.Range('E36').Value = '  +1.48%  '
$ws.Range('E37').Value = '  +1.32%  '
$ws.Range('E38').Value = '  +0.89%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  +4.45%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.36'
$ws.Range('E42').Value = '  -3.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '153.36'
$ws.Range('E43').Value = '  +2.69%  '
$ws.Range('E44').Value = '  +2.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.18'
$ws.Range('E45').Value = '  +3.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.606'
$ws.Range('E46').Value = '  +1.20%  '
$ws.Range('E47').Value = '  +0.61%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0240'
$ws.Range('E48').Value = '  +6.32%  '
$ws.Range('B49').Value = 'Hedera'
$ws.Range('C49').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0524'
$ws.Range('E49').Value = '  +1.87%  '
$ws.Range('E50').Value = '  +3.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.79'
$ws.Range('E51').Value = '  +2.56%  '

Write-Host "Applied 78 cell updates"
